$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.510.27"
$ws.Range("E2").Value = "  -0.25%  "

$ws.Range("D3").Value = "1.617.55"
$ws.Range("E3").Value = "  -1.45%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").Value = "'210.84"
$ws.Range("E5").Value = "  -0.92%  "

$ws.Range("E6").Value = "  -1.87%  "

$ws.Range("E7").Value = "  +0.09%  "

$ws.Range("D8").Value = "'22.79"
$ws.Range("E8").Value = "  -1.29%  "

$ws.Range("D9").Value = "'0.262"
$ws.Range("E9").Value = "  +1.89%  "

$ws.Range("D10").Value = "'0.0611"
$ws.Range("E10").Value = "  +0.04%  "

$ws.Range("D11").Value = "'0.0886"
$ws.Range("E11").Value = "  -0.47%  "

$ws.Range("D12").Value = "1.846.09"
$ws.Range("E12").Value = "  -1.37%  "

$ws.Range("D13").Value = "1.619.21"
$ws.Range("E13").Value = "  -1.28%  "

$ws.Range("E14").Value = "  -0.29%  "

$ws.Range("E15").Value = "  -2.19%  "

$ws.Range("D16").Value = "'64.62"
$ws.Range("E16").Value = "  +0.85%  "

$ws.Range("D17").Value = "27.505.86"
$ws.Range("E17").Value = "  -0.27%  "

$ws.Range("D18").Value = "'229.66"
$ws.Range("E18").Value = "  +0.07%  "

$ws.Range("E19").Value = "  -0.80%  "

$ws.Range("D20").Value = "'7.52"
$ws.Range("E20").Value = "  -1.51%  "

$ws.Range("E21").Value = "  +0.00%  "

$ws.Range("D22").Value = "'4.28"
$ws.Range("E22").Value = "  -0.48%  "

$ws.Range("D23").Value = "'9.99"
$ws.Range("E23").Value = "  +0.41%  "

$ws.Range("D24").Value = "'2.09"
$ws.Range("E24").Value = "  +7.25%  "

$ws.Range("D25").Value = "'148.62"
$ws.Range("E25").Value = "  -0.51%  "

$ws.Range("B26").Value = "BinanceUSD"
$ws.Range("C26").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  +0.13%  "

$ws.Range("B27").Value = "Stellar"
$ws.Range("C27").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D27").Value = "'0.111"
$ws.Range("E27").Value = "  -1.50%  "

$ws.Range("B28").Value = "Cosmos"
$ws.Range("C28").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D28").Value = "'6.80"
$ws.Range("E28").Value = "  -2.29%  "

$ws.Range("D29").Value = "'15.55"
$ws.Range("E29").Value = "  -0.41%  "

$ws.Range("E30").Value = "  -0.58%  "

$ws.Range("D32").Value = "'3.26"
$ws.Range("E32").Value = "  -0.97%  "

$ws.Range("D33").Value = "1.440.68"
$ws.Range("E33").Value = "  +0.91%  "

$ws.Range("E34").Value = "  -3.02%  "

$ws.Range("E35").Value = "  -3.29%  "

$ws.Range("E36").Value = "  +0.03%  "

$ws.Range("D37").Value = "'0.939"
$ws.Range("E37").Value = "  +4.30%  "

$ws.Range("E38").Value = "  -2.25%  "

$ws.Range("D39").Value = "'0.0167"
$ws.Range("E39").Value = "  +0.26%  "

$ws.Range("D40").Value = "'0.861"
$ws.Range("E40").Value = "  -2.06%  "

$ws.Range("D41").Value = "'69.18"
$ws.Range("E41").Value = "  +6.28%  "

$ws.Range("E42").Value = "  +0.11%  "

$ws.Range("D43").Value = "'1.00"
$ws.Range("E43").Value = "  -2.40%  "

$ws.Range("E44").Value = "  +0.06%  "

$ws.Range("D45").Value = "'5.39"
$ws.Range("E45").Value = "  -2.23%  "

$ws.Range("D47").Value = "1.756.92"
$ws.Range("E47").Value = "  -1.36%  "

$ws.Range("D48").Value = "'1.68"
$ws.Range("E48").Value = "  -0.54%  "

$ws.Range("D49").Value = "'86.44"
$ws.Range("E49").Value = "  +0.22%  "

$ws.Range("E50").Value = "  -1.35%  "

$ws.Range("D51").Value = "'0.0992"
$ws.Range("E51").Value = "  +0.73%  "
